$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "CONDITION" column before the existing column D (first ACTION
# column), shifting the two ACTION columns (D,E) one place right (-> E,F).
$ws.Columns.Item(4).Insert()

# Header row: the new column D is a third CONDITION column.
$ws.Range("D18").Value = "CONDITION"

# New condition column label/expression for the "Senior Plus" rule column.
$ws.Range("D19").Value = "Senior Plus"

# Existing rule rows (20-25) don't use the new condition -> leave blank.
$ws.Range("D20").Value = ""
$ws.Range("D21").Value = ""
$ws.Range("D22").Value = ""
$ws.Range("D23").Value = ""
$ws.Range("D24").Value = ""
$ws.Range("D25").Value = ""

# New rule row 26 for the new "Senior Plus 01100949" rule.
$ws.Range("A26").Value = "Senior Plus 01100949"
$ws.Range("D26").Value = "Senior Plus 01100949"
